$wb = $excel.ActiveWorkbook

# --- Sheet: Descriptif_numerique ---
$ws1 = $wb.Worksheets.Item("Descriptif_numerique")
$ws1.Range("C3").Value = 516.177
$ws1.Range("F3").Value = 651.634
$ws1.Range("C4").Value = 670.7828
$ws1.Range("F4").Value = 698.2683
$ws1.Range("C7").Value = 302.005
$ws1.Range("F7").Value = 414.8026
$ws1.Range("C8").Value = 856.75
$ws1.Range("F8").Value = 1014.5591
$ws1.Range("C9").Value = 6327.23
$ws1.Range("F9").Value = 6642.919

# --- Sheet: Regression ---
$ws2 = $wb.Worksheets.Item("Regression")
$ws2.Range("B2").Value = 2.349004527793824
$ws2.Range("B3").Value = 0.9994391188903399
$ws2.Range("B4").Value = 1.082739916204542
$ws2.Range("B5").Value = 1.294326892916941

# --- Sheet: Regression_R2 ---
$ws3 = $wb.Worksheets.Item("Regression_R2")
$ws3.Range("A2").Value = 0.99981413114426
